$d = $word.ActiveDocument

$replacements = @(
    @("487×2=", "390×3="),
    @("493×9=", "780×9="),
    @("918×4=", "789×6="),
    @("496×5=", "844×3="),
    @("368×3=", "398×7="),
    @("360×8=", "927×4="),
    @("903×8=", "947×3="),
    @("744×9=", "821×4="),
    @("972×6=", "302×4="),
    @("886×8=", "128×4="),
    @("793×8=", "775×8="),
    @("623×2=", "381×2="),
    @("633×6=", "896×7="),
    @("969×5=", "268×8="),
    @("380×2=", "116×8="),
    @("439×2=", "499×9="),
    @("787×2=", "132×2="),
    @("900×6=", "762×6="),
    @("751×5=", "764×3="),
    @("432×4=", "225×5="),
    @("992×9=", "221×6="),
    @("361×4=", "574×9="),
    @("807×3=", "697×2="),
    @("484×8=", "930×9="),
    @("283×8=", "760×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
